$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.753.29"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.126.58"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +1.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.68"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5280"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4583"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.58"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09130"
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.58"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.133.31"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.879"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.144"
$ws.Range("E15").Value = "  +5.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001179"
$ws.Range("E16").Value = "  +4.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.62"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.417"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.826.22"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.369"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.373.66"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.58"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.79"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.565"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.65"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.211"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1081"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.397"
$ws.Range("E34").Value = "  +3.81%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.63"
$ws.Range("E36").Value = "  +1.37%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02682"
$ws.Range("E37").Value = "  +4.24%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.909"
$ws.Range("E38").Value = "  +7.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06900"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2335"
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.69"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6933"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.267"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.22"
$ws.Range("E44").Value = "  +8.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6497"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("E47").Value = "  +20.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.708"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.262"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07313"
$ws.Range("E51").Value = "  +3.76%  "
